$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7393.5
$ws.Range("I21").Value = 4504.8335
$ws.Range("J21").Value = 16059.5
$ws.Range("K21").Value = 4504.8335
$ws.Range("L21").Value = 16059.5
$ws.Range("M21").Value = -4036.8335
$ws.Range("N21").Value = -16995.5
$ws.Range("H23").Value = 7393.5
$ws.Range("I23").Value = 4504.8335
$ws.Range("J23").Value = 16059.5
$ws.Range("K23").Value = 4504.8335
$ws.Range("L23").Value = 16059.5
$ws.Range("M23").Value = -4270.8335
$ws.Range("N23").Value = -16527.5
$ws.Range("H33").Value = 22737064
$ws.Range("I33").Value = 50001540
$ws.Range("J33").Value = 16666
$ws.Range("K33").Value = 50001540
$ws.Range("L33").Value = 16666
$ws.Range("M33").Value = -50001311
$ws.Range("N33").Value = -17124
$ws.Range("H40").Value = 6534
$ws.Range("J40").Value = 6534
$ws.Range("L40").Value = 6534
$ws.Range("N40").Value = -6884
$ws.Range("H53").Value = 6579
$ws.Range("I53").Value = 834.5
$ws.Range("J53").Value = 16426.715
$ws.Range("K53").Value = 834.5
$ws.Range("L53").Value = 16426.715
$ws.Range("M53").Value = -197.5
$ws.Range("N53").Value = -17700.715
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2828
$ws.Range("H80").Value = 2161.6667
$ws.Range("I80").Value = 1105.7142
$ws.Range("J80").Value = 3640
$ws.Range("K80").Value = 3317.1426
$ws.Range("L80").Value = 10920
$ws.Range("M80").Value = -2319.1426
$ws.Range("N80").Value = -12916
$ws.Range("H83").Value = 2161.6667
$ws.Range("I83").Value = 1105.7142
$ws.Range("J83").Value = 3640
$ws.Range("K83").Value = 9951.427799999999
$ws.Range("L83").Value = 32760
$ws.Range("M83").Value = -4959.427799999999
$ws.Range("N83").Value = -42744

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1698.3334
$ws.Range("I45").Value = 1396.6666
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1396.6666
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1019.6666
$ws.Range("N45").Value = -2754
$ws.Range("H110").Value = 1174.6471
$ws.Range("I110").Value = 1185.5625
$ws.Range("K110").Value = 1185.5625
$ws.Range("M110").Value = 859.4375
$ws.Range("H132").Value = 2715235.2
$ws.Range("I132").Value = 4554.48
$ws.Range("J132").Value = 8362486.5
$ws.Range("K132").Value = 13663.44
$ws.Range("L132").Value = 25087459.5
$ws.Range("M132").Value = -11133.44
$ws.Range("N132").Value = -25092519.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 642.4286
$ws.Range("I22").Value = 699.4
$ws.Range("K22").Value = 699.4
$ws.Range("M22").Value = -526.4
$ws.Range("H44").Value = 35332.918
$ws.Range("J44").Value = 35332.918
$ws.Range("L44").Value = 35332.918
$ws.Range("N44").Value = -36326.918
$ws.Range("H94").Value = 3034.1667
$ws.Range("I94").Value = 2404
$ws.Range("J94").Value = 6185
$ws.Range("K94").Value = 2404
$ws.Range("L94").Value = 6185
$ws.Range("M94").Value = -1953
$ws.Range("N94").Value = -7087

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4854
$ws.Range("I16").Value = 1078.5555
$ws.Range("K16").Value = 1078.5555
$ws.Range("M16").Value = -791.5554999999999
$ws.Range("H31").Value = 44956
$ws.Range("I31").Value = 27083.166
$ws.Range("K31").Value = 27083.166
$ws.Range("M31").Value = -26788.166
$ws.Range("H34").Value = 44956
$ws.Range("I34").Value = 27083.166
$ws.Range("K34").Value = 27083.166
$ws.Range("M34").Value = -26881.166
$ws.Range("H57").Value = 250000
$ws.Range("J57").Value = 250000
$ws.Range("L57").Value = 250000
$ws.Range("N57").Value = -251508
$ws.Range("H107").Value = 3274.52
$ws.Range("I107").Value = 1441.9375
$ws.Range("J107").Value = 6532.4443
$ws.Range("K107").Value = 1441.9375
$ws.Range("L107").Value = 6532.4443
$ws.Range("M107").Value = 478.0625
$ws.Range("N107").Value = -10372.4443
$ws.Range("H113").Value = 4854
$ws.Range("I113").Value = 1078.5555
$ws.Range("K113").Value = 1078.5555
$ws.Range("M113").Value = 1091.4445
$ws.Range("H120").Value = 33333.332
$ws.Range("I120").Value = 25000
$ws.Range("J120").Value = 50000
$ws.Range("K120").Value = 25000
$ws.Range("L120").Value = 50000
$ws.Range("M120").Value = -21371
$ws.Range("N120").Value = -57258
$ws.Range("H122").Value = 7553.6
$ws.Range("I122").Value = 3859.8
$ws.Range("K122").Value = 11579.4
$ws.Range("M122").Value = -9129.400000000001
$ws.Range("H132").Value = 8069.3184
$ws.Range("I132").Value = 2071
$ws.Range("J132").Value = 16733.555
$ws.Range("K132").Value = 6213
$ws.Range("L132").Value = 50200.665
$ws.Range("M132").Value = -3683
$ws.Range("N132").Value = -55260.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.458336
$ws.Range("I2").Value = 83.882355
$ws.Range("K2").Value = 503.29413
$ws.Range("M2").Value = -390.29413
$ws.Range("H5").Value = 6099536
$ws.Range("I5").Value = 2899
$ws.Range("K5").Value = 8697
$ws.Range("M5").Value = -8585
$ws.Range("H86").Value = 790.8182
$ws.Range("I86").Value = 769.9
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2309.7
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1123.7
$ws.Range("N86").Value = -5372
$ws.Range("H89").Value = 790.8182
$ws.Range("I89").Value = 769.9
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6929.099999999999
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -1001.099999999999
$ws.Range("N89").Value = -20856
$ws.Range("H135").Value = 6099536
$ws.Range("I135").Value = 2899
$ws.Range("K135").Value = 26091
$ws.Range("M135").Value = -23556
$ws.Range("H137").Value = 2183
$ws.Range("J137").Value = 2183
$ws.Range("L137").Value = 6549
$ws.Range("N137").Value = -16749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16260.182
$ws.Range("I70").Value = 20327.818
$ws.Range("J70").Value = 12192.546
$ws.Range("K70").Value = 20327.818
$ws.Range("L70").Value = 12192.546
$ws.Range("M70").Value = -20057.818
$ws.Range("N70").Value = -12732.546
$ws.Range("H73").Value = 16260.182
$ws.Range("I73").Value = 20327.818
$ws.Range("J73").Value = 12192.546
$ws.Range("K73").Value = 20327.818
$ws.Range("L73").Value = 12192.546
$ws.Range("M73").Value = -19391.818
$ws.Range("N73").Value = -14064.546
$ws.Range("H80").Value = 32440.6
$ws.Range("J80").Value = 41767.5
$ws.Range("L80").Value = 41767.5
$ws.Range("N80").Value = -43763.5
$ws.Range("H83").Value = 32440.6
$ws.Range("J83").Value = 41767.5
$ws.Range("L83").Value = 208837.5
$ws.Range("N83").Value = -218821.5
$ws.Range("H132").Value = 9319.883
$ws.Range("I132").Value = 5605.9644
$ws.Range("K132").Value = 16817.8932
$ws.Range("M132").Value = -14287.8932

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2524.25
$ws.Range("I46").Value = 1639.4
$ws.Range("J46").Value = 3999
$ws.Range("K46").Value = 1639.4
$ws.Range("L46").Value = 3999
$ws.Range("M46").Value = -1451.4
$ws.Range("N46").Value = -4375
$ws.Range("H61").Value = 2812.303
$ws.Range("I61").Value = 1811.5
$ws.Range("K61").Value = 1811.5
$ws.Range("M61").Value = -1609.5
$ws.Range("H113").Value = 2812.303
$ws.Range("I113").Value = 1811.5
$ws.Range("K113").Value = 1811.5
$ws.Range("M113").Value = 358.5
$ws.Range("H132").Value = 1495356.8
$ws.Range("I132").Value = 7550.533
$ws.Range("J132").Value = 3355114.5
$ws.Range("K132").Value = 22651.599
$ws.Range("L132").Value = 10065343.5
$ws.Range("M132").Value = -20121.599
$ws.Range("N132").Value = -10070403.5
$ws.Range("H136").Value = 11787.639
$ws.Range("I136").Value = 11295.608
$ws.Range("K136").Value = 33886.824
$ws.Range("M136").Value = -31336.824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1998.5714
$ws.Range("I81").Value = 1117.6
$ws.Range("J81").Value = 4201
$ws.Range("K81").Value = 2235.2
$ws.Range("L81").Value = 8402
$ws.Range("M81").Value = -1174.2
$ws.Range("N81").Value = -10524
$ws.Range("H84").Value = 1998.5714
$ws.Range("I84").Value = 1117.6
$ws.Range("J84").Value = 4201
$ws.Range("K84").Value = 11176
$ws.Range("L84").Value = 42010
$ws.Range("M84").Value = -5872
$ws.Range("N84").Value = -52618
$ws.Range("H113").Value = 3013.375
$ws.Range("I113").Value = 3643.7646
$ws.Range("J113").Value = 1482.4286
$ws.Range("K113").Value = 10931.2938
$ws.Range("L113").Value = 4447.2858
$ws.Range("M113").Value = -8761.293799999999
$ws.Range("N113").Value = -8787.2858
